$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 491, shifting existing rows 491-573 down to 492-574
$ws.Rows.Item(491).Insert()

# Populate the newly inserted row 491 with the new record's data
$ws.Cells.Item(491, 1).Value = 5
$ws.Cells.Item(491, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(491, 3).Value = "Maule"
$ws.Cells.Item(491, 4).Value = 44951
$ws.Cells.Item(491, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(491, 5).Value = 7
$ws.Cells.Item(491, 6).Value = 100112043
$ws.Cells.Item(491, 7).Value = "Pepino ensalada"
$ws.Cells.Item(491, 8).Value = "Sin especificar"
$ws.Cells.Item(491, 9).Value = "Primera"
$ws.Cells.Item(491, 10).Value = 450
$ws.Cells.Item(491, 11).Value = 5000
$ws.Cells.Item(491, 12).Value = 6000
$ws.Cells.Item(491, 13).Value = 5556
$ws.Cells.Item(491, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(491, 15).Value = "Región del Maule"
$ws.Cells.Item(491, 16).Value = 69
$ws.Cells.Item(491, 17).Value = 80
$ws.Cells.Item(491, 18).Value = "Hortaliza"
